$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row needs to be inserted right before the
# existing row 425 ("Femacal de La Calera" - Espinaca records), which
# pushes all subsequent data rows (425-554) down by one (to 426-555).
# The simplest reliable way to do that while keeping every other column
# (Mercado, Region, Categoria, Calidad, Unidad, Origen, Clasificacion, ...)
# identical to its neighbour is to copy row 425 and insert the copy above
# itself, then overwrite just the handful of cells that actually carry
# new data (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg).

$ws.Rows("425:425").Copy()
$ws.Rows("425:425").Insert()

$ws.Range("D425").Value = 45093
$ws.Range("J425").Value = 125
$ws.Range("K425").Value = 5000
$ws.Range("L425").Value = 5500
$ws.Range("M425").Value = 5240
$ws.Range("P425").Value = 1747
